$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the trailing "/12" from the student IDs in column B for rows 71-97
for ($r = 71; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $old = $cell.Value2
    if ($old -like "*/12") {
        $new = $old -replace "/12$", ""
        $cell.Value = $new
    }
}

# Widen column C (name column) to fit the content
$ws.Columns.Item(3).ColumnWidth = 22.4

# Update the active selection / scroll position
$ws.Range("B73").Select() | Out-Null
